# "Add files via upload" / "added more structures data"
#
# Inserts a new row (new row 50, "tensile_stress") ahead of the existing
# "max_stress" row, renames that existing row to "compressive_stress"
# (keeping its existing value/unit), bumps "max_deflection"'s value from
# 3 to 4, and extends the Table1 listobject / AutoFilter / _FilterDatabase
# defined name / conditional formatting down one extra row to match the
# new A1:H66 extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert one blank row at sheet row 50. Everything that was
#    at row 50 (max_stress) through row 65 (max_allowable_wing_span) shifts
#    down by one, to rows 51-66.
$ws.Rows.Item(50).EntireRow.Insert()

# 2) New row 50: tensile_stress / wing stress tensile max.
#    Variable names first (A50, A51) then descriptions (B51, B50) so the
#    shared-string table append order matches the target workbook.
$ws.Range("A50").Value = "tensile_stress"
$ws.Range("A51").Value = "compressive_stress"
$ws.Range("B51").Value = "wing stress compressive max"
$ws.Range("B50").Value = "wing stress tensile max"

$ws.Range("C50").Formula = "=2280000000"
$ws.Range("C50").Font.Color = 0
$ws.Range("C50").Interior.Color = 14277081

$ws.Range("F50").Value = $true

# 3) Row 51 used to be "max_stress" / "wing stress max" with C51=1140000000
#    (Pa) - keep that value/unit/style, it has just been renamed.
$ws.Range("C51").Font.Color = 0

# 4) Row 52 ("max_deflection") value changes from 3 to 4; unit (m) and
#    boolean flag were already carried down by the row insert.
$ws.Range("C52").Value = 4

# 5) Resize the table / AutoFilter to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H66"))

# 6) Update the _xlnm._FilterDatabase defined name to match.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$66"
    }
}

# 7) Extend the conditional formatting range (E2:H65 -> E2:H66).
$fc = $ws.Range("E2").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("E2:H66"))
